# "New Formula to create VM source table"
# Add a new row (row 5) to the Library_Formula sheet, duplicating the pattern
# of the existing CREATE/MODIFY / Utils rows, with a new Formula Name
# "createTableforVM".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# --- Values -------------------------------------------------------------
$ws.Range("A5").Value2 = $ws.Range("A4").Value2
$ws.Range("B5").Value2 = $ws.Range("B4").Value2
$ws.Range("C5").Value2 = "createTableforVM"
$ws.Range("E5").Value2 = $ws.Range("E4").Value2
$ws.Range("F5").Value2 = $ws.Range("F4").Value2

# --- Formatting (mirror the formatting used on row 4) -------------------
$cols = @("A", "B", "C", "E", "F")
foreach ($col in $cols) {
    $src = $ws.Range($col + "4")
    $dst = $ws.Range($col + "5")
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Italic = $src.Font.Italic
    $dst.Font.Color = $src.Font.Color
}

# --- Selection, matching the post-edit cursor position -------------------
$ws.Range("D7").Select()
